$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 28
$ws.Range("H28").Value = 80276.84
$ws.Range("I28").Value = 585.5714
$ws.Range("K28").Value = 585.5714
$ws.Range("M28").Value = -100.5714

# ALC row 43
$ws.Range("H43").Value = 6768.613
$ws.Range("J43").Value = 6447.9165
$ws.Range("L43").Value = 6447.9165
$ws.Range("N43").Value = -6585.9165

# ALC row 45
$ws.Range("H45").Value = 6149.5
$ws.Range("J45").Value = 4900
$ws.Range("L45").Value = 14700
$ws.Range("N45").Value = -15084

# ALC row 116
$ws.Range("H116").Value = 8304
$ws.Range("I116").Value = 7130.3335
$ws.Range("K116").Value = 7130.3335
$ws.Range("M116").Value = -3688.3335

# ALC row 132
$ws.Range("H132").Value = 4545.6
$ws.Range("I132").Value = 4475.6924
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 13427.0772
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -10897.0772
$ws.Range("N132").Value = -20060

# ALC row 137
$ws.Range("H137").Value = 2971.7144
$ws.Range("I137").Value = 1732.3077
$ws.Range("K137").Value = 5196.9231
$ws.Range("M137").Value = -2646.9231

# ALC row 138
$ws.Range("H138").Value = 6221.56
$ws.Range("I138").Value = 5233.909
$ws.Range("J138").Value = 6997.5713
$ws.Range("K138").Value = 15701.727
$ws.Range("L138").Value = 20992.7139
$ws.Range("M138").Value = -10561.727
$ws.Range("N138").Value = -31272.7139

$ws = $wb.Worksheets.Item("ARM")
# ARM row 45
$ws.Range("H45").Value = 2742.625
$ws.Range("J45").Value = 3600
$ws.Range("L45").Value = 3600
$ws.Range("N45").Value = -4354

# ARM row 102
$ws.Range("H102").Value = 2048.8333
$ws.Range("I102").Value = 1406.75
$ws.Range("K102").Value = 1406.75
$ws.Range("M102").Value = 215.25

# ARM row 112
$ws.Range("H112").Value = 189998.5
$ws.Range("J112").Value = 189998.5
$ws.Range("L112").Value = 189998.5
$ws.Range("N112").Value = -192952.5

# ARM row 122
$ws.Range("H122").Value = 3440.125
$ws.Range("I122").Value = 3464.1052
$ws.Range("K122").Value = 10392.3156
$ws.Range("M122").Value = -7942.3156

# ARM row 132
$ws.Range("H132").Value = 8659.799999999999
$ws.Range("I132").Value = 6085.4287
$ws.Range("K132").Value = 18256.2861
$ws.Range("M132").Value = -15726.2861

$ws = $wb.Worksheets.Item("BSM")
# BSM row 20
$ws.Range("H20").Value = 2230.2632
$ws.Range("I20").Value = 1994.4706
$ws.Range("K20").Value = 1994.4706
$ws.Range("M20").Value = -1747.4706

# BSM row 86
$ws.Range("H86").Value = 1301.25
$ws.Range("I86").Value = 1569.25
$ws.Range("J86").Value = 497.25
$ws.Range("K86").Value = 1569.25
$ws.Range("L86").Value = 497.25
$ws.Range("M86").Value = -446.25
$ws.Range("N86").Value = -2743.25

# BSM row 89
$ws.Range("H89").Value = 1301.25
$ws.Range("I89").Value = 1569.25
$ws.Range("J89").Value = 497.25
$ws.Range("K89").Value = 7846.25
$ws.Range("L89").Value = 2486.25
$ws.Range("M89").Value = -2230.25
$ws.Range("N89").Value = -13718.25

# BSM row 94
$ws.Range("H94").Value = 1299.8889
$ws.Range("I94").Value = 1087.375
$ws.Range("K94").Value = 1087.375
$ws.Range("M94").Value = -636.375

# BSM row 134
$ws.Range("H134").Value = 3782.818
$ws.Range("I134").Value = 2067.889
$ws.Range("K134").Value = 6203.667
$ws.Range("M134").Value = -3668.667

$ws = $wb.Worksheets.Item("CRP")
# CRP row 99
$ws.Range("H99").Value = 2453.349
$ws.Range("I99").Value = 2285.1428
$ws.Range("J99").Value = 3042.0715
$ws.Range("K99").Value = 2285.1428
$ws.Range("L99").Value = 3042.0715
$ws.Range("M99").Value = -787.1428000000001
$ws.Range("N99").Value = -6038.0715

# CRP row 126
$ws.Range("H126").Value = 2453.349
$ws.Range("I126").Value = 2285.1428
$ws.Range("J126").Value = 3042.0715
$ws.Range("K126").Value = 6855.428400000001
$ws.Range("L126").Value = 9126.2145
$ws.Range("M126").Value = -4385.428400000001
$ws.Range("N126").Value = -14066.2145

# CRP row 132
$ws.Range("H132").Value = 5606.0625
$ws.Range("I132").Value = 3292.5
$ws.Range("K132").Value = 9877.5
$ws.Range("M132").Value = -7347.5

$ws = $wb.Worksheets.Item("CUL")
# CUL row 63
$ws.Range("H63").Value = 18003.25
$ws.Range("J63").Value = 7337.6665
$ws.Range("L63").Value = 22012.9995
$ws.Range("N63").Value = -23510.9995

# CUL row 66
$ws.Range("H66").Value = 18003.25
$ws.Range("J66").Value = 7337.6665
$ws.Range("L66").Value = 66038.9985
$ws.Range("N66").Value = -73526.9985

# CUL row 107
$ws.Range("H107").Value = 278.2
$ws.Range("I107").Value = 155
$ws.Range("J107").Value = 309
$ws.Range("K107").Value = 465
$ws.Range("L107").Value = 927
$ws.Range("M107").Value = 1455
$ws.Range("N107").Value = -4767

# CUL row 140
$ws.Range("H140").Value = 1725.5
$ws.Range("I140").Value = 1423.7693
$ws.Range("J140").Value = 3033
$ws.Range("K140").Value = 4271.3079
$ws.Range("L140").Value = 9099
$ws.Range("M140").Value = 908.6921000000002
$ws.Range("N140").Value = -19459

$ws = $wb.Worksheets.Item("GSM")
# GSM row 4
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

# GSM row 122
$ws.Range("H122").Value = 5267.3237
$ws.Range("I122").Value = 4039.8462
$ws.Range("J122").Value = 9256.625
$ws.Range("K122").Value = 12119.5386
$ws.Range("L122").Value = 27769.875
$ws.Range("M122").Value = -9669.5386
$ws.Range("N122").Value = -32669.875

# GSM row 132
$ws.Range("H132").Value = 13039.2
$ws.Range("I132").Value = 9732.333000000001
$ws.Range("J132").Value = 17999.5
$ws.Range("K132").Value = 29196.999
$ws.Range("L132").Value = 53998.5
$ws.Range("M132").Value = -26666.999
$ws.Range("N132").Value = -59058.5

$ws = $wb.Worksheets.Item("LTW")
# LTW row 7
$ws.Range("H7").Value = 60571.79
$ws.Range("I7").Value = 90441.086
$ws.Range("J7").Value = 9367.286
$ws.Range("K7").Value = 90441.086
$ws.Range("L7").Value = 9367.286
$ws.Range("M7").Value = -90329.086
$ws.Range("N7").Value = -9591.286

# LTW row 93
$ws.Range("H93").Value = 3069.8572
$ws.Range("I93").Value = 3700
$ws.Range("K93").Value = 3700
$ws.Range("M93").Value = -2452

# LTW row 122
$ws.Range("H122").Value = 7590.5386
$ws.Range("I122").Value = 7219.6665
$ws.Range("J122").Value = 8425
$ws.Range("K122").Value = 21658.9995
$ws.Range("L122").Value = 25275
$ws.Range("M122").Value = -19208.9995
$ws.Range("N122").Value = -30175

# LTW row 126
$ws.Range("H126").Value = 60571.79
$ws.Range("I126").Value = 90441.086
$ws.Range("J126").Value = 9367.286
$ws.Range("K126").Value = 271323.258
$ws.Range("L126").Value = 28101.858
$ws.Range("M126").Value = -268853.258
$ws.Range("N126").Value = -33041.858

# LTW row 132
$ws.Range("H132").Value = 6314.564
$ws.Range("I132").Value = 5899.2705
$ws.Range("K132").Value = 17697.8115
$ws.Range("M132").Value = -15167.8115

$ws = $wb.Worksheets.Item("WVR")
# WVR row 74
$ws.Range("H74").Value = 13166.333
$ws.Range("I74").Value = 12999.5
$ws.Range("K74").Value = 12999.5
$ws.Range("M74").Value = -12063.5

# WVR row 77
$ws.Range("H77").Value = 13166.333
$ws.Range("I77").Value = 12999.5
$ws.Range("K77").Value = 38998.5
$ws.Range("M77").Value = -34318.5

# WVR row 107
$ws.Range("H107").Value = 1429441
$ws.Range("I107").Value = 1818838.9
$ws.Range("K107").Value = 5456516.699999999
$ws.Range("M107").Value = -5454596.699999999

# WVR row 122
$ws.Range("H122").Value = 4264.8335
$ws.Range("I122").Value = 4400
$ws.Range("J122").Value = 3589
$ws.Range("K122").Value = 13200
$ws.Range("L122").Value = 10767
$ws.Range("M122").Value = -10750
$ws.Range("N122").Value = -15667

# WVR row 132
$ws.Range("H132").Value = 3982
$ws.Range("I132").Value = 2256.8572
$ws.Range("J132").Value = 6397.2
$ws.Range("K132").Value = 6770.571599999999
$ws.Range("L132").Value = 19191.6
$ws.Range("M132").Value = -4240.571599999999
$ws.Range("N132").Value = -24251.6

# WVR row 138
$ws.Range("H138").Value = 80139.664
$ws.Range("J138").Value = 80139.664
$ws.Range("L138").Value = 80139.664
$ws.Range("N138").Value = -90419.664
